$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.619.27'
$ws.Range("E2").Value = '  +1.97%  '
$ws.Range("D3").Value = '2.283.59'
$ws.Range("E3").Value = '  +3.59%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.98'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.627'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.80'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +9.69%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.649'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +5.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.28'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0987'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.97'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.31'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +4.51%  '
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("D15").Value = '2.627.68'
$ws.Range("E15").Value = '  +3.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.06'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +4.26%  '
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").Value = '2.281.00'
$ws.Range("E18").Value = '  +2.80%  '
$ws.Range("D19").Value = '42.569.89'
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("E20").Value = '  +5.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.32'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.27'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.92'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.51%  '
$ws.Range("E24").Value = '  +9.24%  '
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.41'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.43'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.72%  '
$ws.Range("E29").Value = '  -0.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.19'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.86'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.09'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.34'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +8.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.128'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +5.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0816'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.96'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +23.45%  '
$ws.Range("E37").Value = '  +3.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.70'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +15.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.75'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0306'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.89'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +16.80%  '
$ws.Range("E42").Value = '  +5.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.94'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.214'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +9.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.16'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +7.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.93'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.85'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.08%  '
$ws.Range("E48").Value = '  +5.05%  '
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.18'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.79'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.37%  '
